$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B column values (B2:B13) with the new model parameters.
$ws.Range("B2").Value = 1.042309107091355
$ws.Range("B3").Value = -0.116562815773544
$ws.Range("B4").Value = -0.1592239461346836
$ws.Range("B5").Value = -0.1228135972531265
$ws.Range("B6").Value = -0.6428507070351096
$ws.Range("B7").Value = 0.1385793388500075
$ws.Range("B8").Value = 0.0841738758598569
$ws.Range("B9").Value = 0.0001971396443468134
$ws.Range("B10").Value = 0.5797073516238369
$ws.Range("B11").Value = -0.05030276577666537
$ws.Range("B12").Value = 0.008861290341707543
$ws.Range("B13").Value = 0.2767204933162186

# Row 14 (shot_during_regular_play) is removed entirely from the model.
$ws.Rows.Item(14).Delete()
